{"js": "// Apply three textual edits to the John Dehner biography document:\n// 1. \"Have Gun [en dash] Will Travel\"  ->  \"Have Gun[em dash]Will Travel\"\n// 2. \"...Frontier Gentleman series. He also...\" ->\n//    \"...Frontier Gentleman series, which immediately preceded it. He also...\"\n// 3. \"(1970-71)\" -> \"(1970-71)\" and \"(1979-80)\" -> \"(1979-80)\" with the\n//    hyphen replaced by an en dash.\n\nconst EN_DASH = \"\\u2013\"; // \u2013\nconst EM_DASH = \"\\u2014\"; // \u2014\n\nconst body = context.document.body;\n\n// --- Edit 1: \"Have Gun <space><en dash><space> Will Travel\" -> \"Have Gun<em dash>Will Travel\" ---\n// Only the dash run (already italic, between \"Have Gun\" and \"Will Travel\") is touched,\n// so the italic formatting on the surrounding runs is left intact.\nconst dashResults = body.search(\" \" + EN_DASH + \" \", { matchCase: true });\ndashResults.load(\"items\");\nawait context.sync();\nif (dashResults.items.length > 0) {\n  dashResults.items[0].insertText(EM_DASH, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Edit 2: insert \", which immediately preceded it\" before the period after \"series\" ---\n// Search only within the (non-italic) run(s) that hold \" series. He also lent\"\n// so the italic \"Frontier Gentleman\" run immediately before it is left untouched.\nconst seriesResults = body.search(\" series. He also lent\", { matchCase: true });\nseriesResults.load(\"items\");\nawait context.sync();\nif (seriesResults.items.length > 0) {\n  seriesResults.items[0].insertText(\n    \" series, which immediately preceded it. He also lent\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// --- Edit 3: hyphen -> en dash in the two television year ranges ---\nconst range1970 = body.search(\"1970-71\", { matchCase: true });\nrange1970.load(\"items\");\nawait context.sync();\nif (range1970.items.length > 0) {\n  range1970.items[0].insertText(\"1970\" + EN_DASH + \"71\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst range1979 = body.search(\"1979-80) as Marshall Edge Troy.\", { matchCase: true });\nrange1979.load(\"items\");\nawait context.sync();\nif (range1979.items.length > 0) {\n  range1979.items[0].insertText(\n    \"1979\" + EN_DASH + \"80) as Marshall Edge Troy.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Apply three textual edits to the John Dehner biography document:\n# 1. \"Have Gun \u2013 Will Travel\"  ->  \"Have Gun\u2014Will Travel\"   (en dash + spaces -> em dash)\n# 2. \"...Frontier Gentleman series. He also...\" -> \"...Frontier Gentleman series, which immediately preceded it. He also...\"\n# 3. \"(1970-71)\" -> \"(1970\u201371)\" and \"(1979-80)\" -> \"(1979\u201380)\"  (hyphen -> en dash)\n\n$d = $word.ActiveDocument\n\n$enDash = [string][char]0x2013\n$emDash = [string][char]0x2014\n\n# --- Edit 1: em dash between \"Have Gun\" and \"Will Travel\" ---\n$find1 = $d.Content.Find\n$find1.Text = \" \" + $enDash + \" \"\n$find1.Replacement.Text = $emDash\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# --- Edit 2: insert \", which immediately preceded it\" before the period after \"series\" ---\n# Search only within the (non-italic) run(s) that hold \" series. He also lent\" so the\n# italic \"Frontier Gentleman\" run immediately before it is left untouched.\n$find2 = $d.Content.Find\n$find2.Text = \" series. He also lent\"\n$find2.Replacement.Text = \" series, which immediately preceded it. He also lent\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# --- Edit 3: hyphen -> en dash in the two year ranges ---\n$find3 = $d.Content.Find\n$find3.Text = \"1970-71\"\n$find3.Replacement.Text = \"1970\" + $enDash + \"71\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2) | Out-Null\n\n$find4 = $d.Content.Find\n$find4.Text = \"1979-80) as Marshall Edge Troy.\"\n$find4.Replacement.Text = \"1979\" + $enDash + \"80) as Marshall Edge Troy.\"\n$find4.Execute($find4.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find4.Replacement.Text, 2) | Out-Null\n"}
